$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (row 2 through row 12): change date value 45184 -> 45185
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
